$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update student IDs in column B (rows 2-13): 2000XXXX -> 2001XXXX
for ($r = 2; $r -le 13; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $old = $cell.Value2
    $cell.Value2 = $old + 10000
}

# Update the active selection from C16 to C15
$ws.Range("C15").Select()
